$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for first file updated
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 07:03:21"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime" for first file updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 07:03:16"
$wsZhCn.Range("K2").Value = "2016-08-26 07:03:32"

# de-de sheet: "Correspond Handoff Datetime" (same value as Overview!G2) and "Correspond Handback DateTime" updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-26 07:03:21"
$wsDeDe.Range("K2").Value = "2016-08-26 07:03:38"
